$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (old D:K shifts right to F:M)
$ws.Columns("D:E").Insert()

# New D,E columns currently have no explicit number format; copy the format
# from column F (the shifted original column D) which retains the correct style.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Header / period-ending date rows: two new quarters inserted at D,E ----
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373

# ---- Data rows: new D,E,F values (old D is replaced; old E:K shifted to G:M already) ----
$dataMap = @{}
$dataMap[8] = @(879100, 909800, 1621800)
$dataMap[9] = @(826800, 877700, 1513900)
$dataMap[10] = @(52300, 32100, 107900)
$dataMap[11] = @($null, $null, $null)
$dataMap[12] = @("NA", "NA", "NA")
$dataMap[13] = @(0, 0, 0)
$dataMap[14] = @(11000, 2100, 1400)
$dataMap[15] = @(13600, 13200, 25800)
$dataMap[16] = @($null, $null, $null)
$dataMap[17] = @(862900, 904900, 1565200)
$dataMap[18] = @(16200, 4900, 56600)
$dataMap[19] = @($null, $null, $null)
$dataMap[20] = @(7600, 0, -5800)
$dataMap[21] = @(37400, 18200, 76600)
$dataMap[22] = @(10400, 10400, 18900)
$dataMap[23] = @(13300, -5500, 31900)
$dataMap[24] = @(-600, 400, 500)
$dataMap[25] = @(0, 0, 0)
$dataMap[26] = @(13900, -5800, 31400)
$dataMap[27] = @(13700, -5800, 31000)
$dataMap[28] = @(0, 0, 0)
$dataMap[29] = @(0, 0, 0)
$dataMap[30] = @(0, 0, 0)
$dataMap[31] = @(0, 0, 0)
$dataMap[32] = @(-7600, 0, 5800)
$dataMap[33] = @(13700, -5800, 31000)
$dataMap[34] = @(0, 0, 0)
$dataMap[35] = @(13700, -5800, 31000)
$dataMap[39] = @($null, $null, $null)
$dataMap[40] = @($null, $null, $null)
$dataMap[41] = @(75100, 87700, 82700)
$dataMap[42] = @(0, 0, 0)
$dataMap[43] = @(160300, 133000, 133500)
$dataMap[44] = @(322100, 358600, 333700)
$dataMap[45] = @(29100, 11000, 87700)
$dataMap[46] = @(586600, 590300, 637600)
$dataMap[47] = @(136700, 131500, 130400)
$dataMap[48] = @(538300, 482800, 477600)
$dataMap[49] = @(177300, 178000, 179100)
$dataMap[50] = @(0, 0, 0)
$dataMap[51] = @(0, 0, 0)
$dataMap[52] = @(21900, 23900, 26600)
$dataMap[53] = @(0, 0, 0)
$dataMap[54] = @(1460700, 1406500, 1451300)
$dataMap[55] = @($null, $null, $null)
$dataMap[56] = @($null, $null, $null)
$dataMap[57] = @(54800, 68600, 52700)
$dataMap[58] = @(0, "NA", "NA")
$dataMap[59] = @(452400, 420000, 479200)
$dataMap[60] = @(507200, 488600, 531900)
$dataMap[61] = @(398700, 395300, 394400)
$dataMap[62] = @(42500, 45200, 43300)
$dataMap[63] = @(0, 0, 0)
$dataMap[64] = @(0, 0, 0)
$dataMap[65] = @(0, 0, 0)
$dataMap[66] = @(948400, 929100, 969700)
$dataMap[67] = @($null, $null, $null)
$dataMap[68] = @(0, 0, 0)
$dataMap[69] = @(0, 0, 0)
$dataMap[70] = @(0, 0, 0)
$dataMap[71] = @(0, 0, 0)
$dataMap[72] = @(-108800, -122600, -116800)
$dataMap[73] = @(0, 0, 0)
$dataMap[74] = @(0, 0, 0)
$dataMap[75] = @(0, 0, 0)
$dataMap[76] = @(512300, 477400, 481600)
$dataMap[77] = @(0, 0, 0)
$dataMap[81] = @(13700, -5800, 31000)
$dataMap[82] = @($null, $null, $null)
$dataMap[83] = @(13600, 13200, 25800)
$dataMap[84] = @(0, 0, 0)
$dataMap[85] = @(0, 0, 0)
$dataMap[86] = @(0, 0, 0)
$dataMap[87] = @(0, 0, 0)
$dataMap[88] = @(0, 0, 0)
$dataMap[89] = @(38700, 20600, 31300)
$dataMap[90] = @($null, $null, $null)
$dataMap[91] = @(-18200, -12500, -17700)
$dataMap[92] = @(0, 0, 0)
$dataMap[93] = @(0, 0, 0)
$dataMap[94] = @(-72100, -12500, -91200)
$dataMap[95] = @($null, $null, $null)
$dataMap[96] = @(0, 0, 0)
$dataMap[97] = @(0, 0, 0)
$dataMap[98] = @(0, 0, 0)
$dataMap[99] = @(0, 0, 0)
$dataMap[100] = @(20700, -3000, 24200)
$dataMap[101] = @(0, 0, 0)
$dataMap[102] = @(-12700, 5000, -35600)

foreach ($row in $dataMap.Keys) {
    $vals = $dataMap[$row]
    if ($vals[0] -ne $null) { $ws.Cells.Item($row, 4).Value2 = $vals[0] }
    if ($vals[1] -ne $null) { $ws.Cells.Item($row, 5).Value2 = $vals[1] }
    if ($vals[2] -ne $null) { $ws.Cells.Item($row, 6).Value2 = $vals[2] }
}

Write-Output "edit complete"
